$d = $word.ActiveDocument

$old1 = 'Appendix 3: SWIFT Qualitative Focus Group Guide: Participants'
$new1 = 'ISihlomelo sesi-3: ISikhokelo soMgangatho seQela eliGxilileyo le-SWIFT: Abathathi-nxaxheba'
$result1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
Write-Host "Replace 1: $result1"

$old2 = 'Briefing:'
$new2 = 'Ingxelo:'
$result2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
Write-Host "Replace 2: $result2"

$old3 = 'Thank the interviewees for making the time to attend the focus group discussion. '
$new3 = 'Bulela abantu ebebesenziwa udliwano-ndlebe ngokwenza ixesha lokuzimasa ingxoxo yeqela ekugxilwe kulo. '
$result3 = $d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)
Write-Host "Replace 3: $result3"

$old4 = 'Thank you for participating in this focus group. You would have all had a telephone conversation with one of the members of our research team, where they went over all the information for this focus group, and you agreed to be part of the study. For some of you, this might have been a while back and so we will cover some of this information briefly as a reminder. '
$new4 = 'Enkosi kwakhona ngokuthatha inxaxheba kweliqela eligxilileyo. Nonke ngenibe nencoko yomnxeba nelinye lamalungu eqela lethu lophando, apho baye bajongisisa lonke ulwazi lweli qela ekugxilwe kulo, kwaye uvumile ukuba yinxalenye yolu phononongo. Kwabanye benu, oku kusenokuba bekusemva kwexesha kwaye ke siza kugubungela ezinye zolu lwazi ngokufutshane njengesikhumbuzo. '
$result4 = $d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2)
Write-Host "Replace 4: $result4"

$old5 = 'This discussion will take about 1-1.5 hours. There are no right or wrong answers. We want to learn about your experiences and what you liked or didn''t like about the parenting programme on WhatsApp, so that we can make it better for other parents and caregivers. The programme is called ParentText. The other members of this focus group are also parents and caregivers like you. We will need to record this discussion, with your permission, so that we can listen to your contributions at a later stage. However, all your answers will be completely confidential and will only be viewed by the research team.'
$new5 = 'Lengxoxo izakuthatha malunga neyure ukuya kwiyure enamashumi amahlanu. Akukho zimpendulo zilungileyo okanye zingalunganga. Sifuna ukufunda malunga namava akho kunye nento oyithandileyo okanye ongakhange uyithande malunga nenkqubo yobuzali kuWhatsApp, ukuze siyenze ukuba ibengcono kwabanye abazali kunye nabanonopheli. Lenkqubo ibizwa ngokuba yi-ParentText. Amanye amalungu eliqela eligxilileyo nawo angabazali kunye nabanonopheli kanje ngawe. We will need to record this discussion, with your permission, so that we can listen to your contributions at a later stage. Nangona kunjalo, zonke iimpendulo zakho ziya kuba yimfihlo ngokupheleleyo kwaye ziya kujongwa kuphela liqela lophando.'
$result5 = $d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2)
Write-Host "Replace 5: $result5"

$old6 = 'We will not be introducing ourselves in the group. This is not something we have forgotten to do, it’s because we want your names and who you are to remain private because this session is being recorded. We will also not use any of your names when we write up everything after this discussion. You have promised that you will respect others in the group and will not discuss what is said by others outside of the group. '
$new6 = 'Asizukuzazisa iziqu zethu kwiqela. This is not something we have forgotten to do, it’s because we want your names and who you are to remain private because this session is being recorded. Kananjalo asizukusebenzisa nawaphi na amagama akho xa sibhala yonke into emva kwale ngxoxo. Uthembisile ukuba uya kubahlonipha abanye eqeleni kwaye awuzukuxoxa ngezinto ezithethwa ngabanye ngaphandle kweqela. '
$result6 = $d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $new6, 2)
Write-Host "Replace 6: $result6"

$old7 = 'We know that during the programme, some sensitive things may have come up, especially in the second module about child safety and perhaps in some of the referrals you chose to access. Because this is a group discussion, we don''t expect you to share any of these sensitive details. Just remember, we only want you to share what you’re comfortable with.'
$new7 = 'Siyazi ukuba ngexesha lenkqubo, ezinye izinto ezinobuzaza zisenokuba ziye zavela, ingakumbi kwimodyuli yesibini emalunga nokhuseleko lwabantwana kwaye mhlawumbi kwezinye zokuthunyelwa okhethe ukufikelela kuzo. Kuba le yingxoxo yeyeqela, asilindelanga ukuba wabelane nangaziphi na kwezi nkcukacha zibuthathaka. Khumbula nje, sifuna kuphela ukuba wabelane ngalento ukhululekileyo ngayo.'
$result7 = $d.Content.Find.Execute($old7, $true, $false, $false, $false, $false, $true, 1, $false, $new7, 2)
Write-Host "Replace 7: $result7"

$old8 = 'Please remember that you are able to leave this discussion at any time if you would like to, or you can choose not to answer any questions that you don’t want to respond to, for any reason. If you decide at a later stage that you would like your contribution to be removed from the study, you can contact the research team by email on your consent form and up until the [date to be determined]. After this point, we would have started to share our findings. Before signing the consent forms, do you have any questions? '
$new8 = 'Nceda ukhumbule ukuba uyakwazi ukuyishiya le ngxoxo nangaliphi na ixesha ukuba uyafuna, okanye ungakhetha ukungaphenduli nayiphi na imibuzo ongafuniyo ukuyiphendula, nangasiphi na isizathu. Ukuba uthatha isgqibo sekumva sokuba ungathanda igalelo lakho lisuswe koluphononongo, ungaqhagamshelana neqela lophando nge-imeyile ikwifomu yakho yemvume kwaye kude kube [ngumhla oza kumiselwa]. Emva kweli nqanaba, ngesele siqalisile ukwabelana ngeziphumo zethu. Ngaphambi kokuba utyikitye ifomu yemvume, ingaba unayo nayiphi na imibuzo? '
$result8 = $d.Content.Find.Execute($old8, $true, $false, $false, $false, $false, $true, 1, $false, $new8, 2)
Write-Host "Replace 8: $result8"

$old9 = 'At this point hand out the consent forms. Tell participants to detach the last page and sign it if they still consent to being part of the study. They can fold the last page in half and hand it back to you as you go around (don’t let them pass it between themselves). Please remind them that all the information they need about the study is on the two pages they still have and that they should keep those safely as it also has contact numbers on it. '
$new9 = 'Kweli nqanaba nikezela ngeefomu zemvume. Tell participants to detach the last page and sign it if they still consent to being part of the study. They can fold the last page in half and hand it back to you as you go around (don’t let them pass it between themselves). Please remind them that all the information they need about the study is on the two pages they still have and that they should keep those safely as it also has contact numbers on it. '
$result9 = $d.Content.Find.Execute($old9, $true, $false, $false, $false, $false, $true, 1, $false, $new9, 2)
Write-Host "Replace 9: $result9"
